# Update "想去人数" (want-to-go count) figures in the F column on both the
# "展览" and "全部类型" sheets (they carry duplicate data).
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 11398
    "F3"  = 10756
    "F6"  = 977
    "F11" = 10566
    "F12" = 4078
    "F16" = 27
    "F18" = 416
    "F19" = 11090
    "F20" = 10844
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
